$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Date column (B) with new test-run timestamps for the
# "Personal_SearchTransaction_Generic_TC" test data rows.
$ws.Range("B2").Value = "Fri Mar 08 01:12:58 EST 2024"
$ws.Range("B3").Value = "Fri Mar 08 01:13:35 EST 2024"
$ws.Range("B5").Value = "Fri Mar 08 01:14:03 EST 2024"
$ws.Range("B6").Value = "Fri Mar 08 01:14:39 EST 2024"
$ws.Range("B7").Value = "Fri Mar 08 01:15:15 EST 2024"
